$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.704.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.86%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.805.49'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.81%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '

$ws.Range("E6").Value = '  +0.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5327'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.35%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3783'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07514'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.117'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.81%  '

$ws.Range("E12").Value = '  +0.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.177'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.363'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.804.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001064'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06449'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.908'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.720.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.101'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.014.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.363'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.102'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1054'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.654'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.07%  '

$ws.Range("E34").Value = '  +2.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2256'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06431'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02311'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.793'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.32%  '

$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.223'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6241'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.397'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5867'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.694'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.951'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.150'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06890'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.78%  '
